$wb = $excel.ActiveWorkbook

# Overview sheet: update "Latest HO Xliff Generate Date" for the first file row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-31 12:49:41"

# zh-cn sheet: update "Correspond Handoff Datetime" and "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-31 12:49:27"
$wsZhCn.Range("K2").Value = "2016-08-31 12:50:34"

# de-de sheet: update "Correspond Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-31 12:50:52"
